$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("J1").Value = "MHC A3 G246"
$ws.Range("J2").Value = "GT"
